$d = $word.ActiveDocument

# The document contains two bullet lists that happen to share wording:
#  1) "Partner - Siege Analytics" > "Data Science & Political Analytics" (unchanged)
#  2) "KEY ACHIEVEMENTS AND IMPACT" > "Impact" (target of this edit)
# So we locate the target paragraphs via the unique "KEY ACHIEVEMENTS AND IMPACT"
# heading's document-order index rather than doing a document-wide text replace.

$i = 0
$startIndex = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Trim() -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $startIndex = $i
    }
}

# Bullets begin two paragraphs after the section heading (heading, then "Impact" sub-heading).
$bullet1 = $d.Paragraphs.Item($startIndex + 2)
$bullet2 = $d.Paragraphs.Item($startIndex + 3)
$bullet3 = $d.Paragraphs.Item($startIndex + 4)
$bullet4 = $d.Paragraphs.Item($startIndex + 5)
$bullet5 = $d.Paragraphs.Item($startIndex + 6)
$bullet6 = $d.Paragraphs.Item($startIndex + 7)

# Rewrite the first four bullets as impact-focused accomplishment statements.
$bullet1.Range.Text = "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
$bullet2.Range.Text = "• `$4.7M savings enabled nonprofit access"
$bullet3.Range.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"
$bullet4.Range.Text = "• 178% accuracy improvement in racial classification algorithms"

# Drop the last two bullets entirely (delete from the end backward so indices stay valid).
$bullet6.Range.Delete()
$bullet5.Range.Delete()
